# The commit removes the document's leading title ("Integrated Capstone
# Project"), the "This Case Study has four check points..." paragraph, the
# blank paragraph that followed it, and the entire "Check Point Topics"
# summary table, leaving the document starting at the paragraph that
# contains the lone tab character (immediately before "Domain:").

$d = $word.ActiveDocument

# Delete the table first. Deleting a Range that spans across a table
# (paragraphs + table together) does not reliably remove the table
# content in this engine, so the table must be removed as its own step.
if ($d.Tables.Count -gt 0) {
    $d.Tables(1).Delete()
}

# After the table is gone, the title paragraph, the "four check points"
# paragraph, and the trailing blank paragraph are the first three
# paragraphs in the document. Remove them by deleting the range that
# covers paragraphs 1-3, leaving paragraph 4 (the tab character) as the
# new start of the body.
$thirdPara = $d.Paragraphs(3)
$deleteEnd = $thirdPara.Range.End
$leadInRange = $d.Range(0, $deleteEnd)
$leadInRange.Delete()
